# Insert a new "Plumb" / Construction row into the Cluster_Keywords table,
# right after the existing "Excav" row (which is the last Construction
# entry), matching the alphabetical sort order used by the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is currently sorted by Cluster Category then Stem, and "Excav"
# (row 14) is the last "Construction" row, with "Justi" (row 15) starting
# the "Corrections" group. Insert a whole worksheet row above the current
# row 15 so the new data lands between them as row 15, shifting the rest
# (including the table) down.
$ws.Rows(15).Insert()

$ws.Range("A15").Value = "Plumb"
$ws.Range("B15").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Range("C15").Value = "Construction"

# The table-row insert rewrites the last row's calculated-column formula
# to the "[@Stem]" shorthand; restore the canonical structured-reference
# form used throughout the rest of the table.
$ws.Range("B73").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"

# Make sure the table (range, autofilter, sort state) fully covers the
# grown data range A1:C73.
$table = $ws.ListObjects.Item("Cluster_Keywords")
$table.Resize($ws.Range("A1:C73"))

# Keep the duplicate-highlighting conditional formats in sync with the
# rows they originally tracked, now shifted down by one.
$ws.Range("A29:A30").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A30:A31"))
$ws.Range("A22:A23").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A23:A24"))
$ws.Range("B2:B72").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B2:B73"))
$ws.Range("A2:A21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A2:A22,A61:A73,A25:A29,A32:A59"))

# Select the cell below the newly entered data, matching where the user's
# cursor ended up after finishing the edit.
$ws.Range("A16").Select()
